# Generate Report for Handback
# Adds a new handback record (4b378967-7a06-41a9-b2c6-b02beb34cc96.md) as
# row 4 to the Overview, zh-cn and de-de report tables.

$wb = $excel.ActiveWorkbook

$fileGuid   = "4b378967-7a06-41a9-b2c6-b02beb34cc96"
$srcRepo    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dec74015c7bd6d5c153ed185be4553477c7e1246/e2e/$fileGuid.md"
$zhcnRepo   = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/5c2bb49d95dc3b03b7c183c0bebb033cb6e48dbb/e2e/$fileGuid.md"
$dedeRepo   = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/221aa6a4a475d14786cf8b08f8f5a4fea4807698/e2e/$fileGuid.md"

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = "$fileGuid.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), $srcRepo, "", "", "e2e\$fileGuid.md") | Out-Null
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G4").Value = "2016-09-09 08:40:43"
$wsOverview.Range("G4").NumberFormat = $dateFmt

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), $srcRepo, "", "", "$fileGuid.md") | Out-Null
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "'True"
$wsZhCn.Range("G4").Value = "$fileGuid.f111c09c0d0ce47bfe9880ebc72de33c0fc98365.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-09-09 08:40:32"
$wsZhCn.Range("H4").NumberFormat = $dateFmt
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I4"), $zhcnRepo, "", "", "$fileGuid.md") | Out-Null
$wsZhCn.Range("J4").Value = "$fileGuid.f111c09c0d0ce47bfe9880ebc72de33c0fc98365.zh-cn.xlf"
$wsZhCn.Range("K4").Value = "2016-09-09 08:41:29"
$wsZhCn.Range("K4").NumberFormat = $dateFmt
$wsZhCn.Range("M4").Value = "'True"
$wsZhCn.Range("O4").Value = "'False"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), $srcRepo, "", "", "$fileGuid.md") | Out-Null
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "'True"
$wsDeDe.Range("G4").Value = "$fileGuid.f111c09c0d0ce47bfe9880ebc72de33c0fc98365.de-de.xlf"
$wsDeDe.Range("H4").Value = "2016-09-09 08:40:43"
$wsDeDe.Range("H4").NumberFormat = $dateFmt
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I4"), $dedeRepo, "", "", "$fileGuid.md") | Out-Null
$wsDeDe.Range("J4").Value = "$fileGuid.f111c09c0d0ce47bfe9880ebc72de33c0fc98365.de-de.xlf"
$wsDeDe.Range("K4").Value = "2016-09-09 08:41:47"
$wsDeDe.Range("K4").NumberFormat = $dateFmt
$wsDeDe.Range("M4").Value = "'True"
$wsDeDe.Range("O4").Value = "'False"
